# Auto-generated script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.943.91"
$ws.Range("E2").Value = "  +4.57%  "
$ws.Range("D3").Value = "'3.516.79"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'593.43"
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("D6").Value = "'168.86"
$ws.Range("E6").Value = "  +5.90%  "
$ws.Range("D8").Value = "'3.516.64"
$ws.Range("E8").Value = "  +2.40%  "
$ws.Range("D9").Value = "'0.576"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("D12").Value = "'0.439"
$ws.Range("E12").Value = "  +3.67%  "
$ws.Range("D13").Value = "'4.126.58"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "'0.135"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "'28.18"
$ws.Range("E15").Value = "  +3.57%  "
$ws.Range("E16").Value = "  +3.28%  "
$ws.Range("D17").Value = "'66.905.99"
$ws.Range("E17").Value = "  +4.42%  "
$ws.Range("D18").Value = "'3.533.95"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "'14.05"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").Value = "'390.81"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "'73.58"
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("E24").Value = "  +9.52%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "'0.532"
$ws.Range("E26").Value = "  +3.18%  "
$ws.Range("D27").Value = "'10.20"
$ws.Range("E27").Value = "  +5.26%  "
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +5.73%  "
$ws.Range("D31").Value = "'1.47"
$ws.Range("E31").Value = "  +5.77%  "
$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'7.48"
$ws.Range("E33").Value = "  +7.17%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'23.59"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "'1.61"
$ws.Range("E36").Value = "  +5.81%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +6.00%  "
$ws.Range("D39").Value = "'1.92"
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("D41").Value = "'4.65"
$ws.Range("E41").Value = "  +6.24%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'26.56"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.831.80"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'6.66"
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").Value = "'43.54"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "'26.30"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'0.0314"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("E48").Value = "  +4.58%  "
$ws.Range("D49").Value = "'351.54"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("D51").Value = "'33.66"
